$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ------------------------------------------------------------------

$titlePara = $d.Paragraphs.Item(1)
$titleText = $titlePara.Range.Text.TrimEnd([char]13, [char]7)

# Anchor the new paragraph next to an existing Normal-style paragraph
# (paragraph 3, "Step right up folks...") so the inserted paragraph
# naturally inherits the "Normal" style without leaving stray rsid
# bookkeeping attributes behind.
$anchor = $d.Paragraphs.Item(3)
$anchor.Range.InsertParagraphBefore()
$metaPara = $d.Paragraphs.Item(3)

$metaText = "Meta description: Discover all about Ankh of Anubis, an online slot game from Play" + [char]0x2019 + "N" + [char]0x2019 + "Go, with an Ancient Egypt theme centered on the god Anubis. Play it free and read our review."

$insertionPoint = $d.Range($metaPara.Range.Start, $metaPara.Range.Start)
$insertionPoint.InsertAfter($metaText)

# Bold just the leading "Meta description" label (16 characters).
$boldRange = $d.Range($metaPara.Range.Start, $metaPara.Range.Start + 16)
$boldRange.Font.Bold = 1

# Move this freshly built paragraph (including its paragraph mark) so
# that it sits right after the first (title) paragraph.
$metaPara = $d.Paragraphs.Item(3)
$cutRange = $d.Range($metaPara.Range.Start, $metaPara.Range.End)
$cutRange.Cut()

$titlePara = $d.Paragraphs.Item(1)
$dest = $d.Range($titlePara.Range.End, $titlePara.Range.End)
$dest.Paste()

# ------------------------------------------------------------------
# 2) Remove the duplicated bold title paragraph near the end of the
#    document, and replace the text of the final (italic) paragraph
#    with the new image-generation prompt, keeping its formatting.
# ------------------------------------------------------------------

# Locate the duplicate title paragraph (bold run with the same text as
# the document title) by content rather than a hard-coded index, so
# the script stays correct even if paragraph numbering shifts.
$dupIndex = -1
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq $titleText) {
        $dupIndex = $i
    }
}

if ($dupIndex -ne -1) {
    $dupTitlePara = $d.Paragraphs.Item($dupIndex)
    $dupTitlePara.Range.Delete()
}

$total = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($total)
$lastTextRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)

$newPromptText = 'Please create a cartoon image for the online slot game "Ankh of Anubis". The image must feature a happy Maya warrior wearing glasses. The style of the image should be cartoonish, with bright colors to catch the player''s attention. The Maya warrior should be depicted in a pose of victory with a thumbs up, as if he has just won a big prize. The background can feature elements of Ancient Egypt, such as pyramids or hieroglyphics. The image should be visually appealing and reflect the exciting and adventurous nature of the game.'

$lastTextRange.Text = $newPromptText
